# Auto-generated edit script: update crypto price/volume columns (D, E)
# for rows 2-51 on Sheet1, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.715.62"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "1.805.56"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.46"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5933"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06825"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.32"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07514"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "1.805.86"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.757"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6232"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").Value = "2.050.76"
$ws.Range("E15").Value = "  -1.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009261"
$ws.Range("E16").Value = "  -6.73%  "
$ws.Range("E17").Value = "  -4.15%  "
$ws.Range("D18").Value = "28.678.99"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.475"
$ws.Range("E19").Value = "  -6.55%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "210.87"
$ws.Range("E21").Value = "  -6.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.50"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.821"
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.96"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.873"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.42"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.429"
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06155"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.785"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.747"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.731"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.062"
$ws.Range("E35").Value = "  -5.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6417"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.499"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.718"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.578"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01697"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").Value = "1.145.98"
$ws.Range("E41").Value = "  -5.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8812"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.008"
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "1.957.18"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.55"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000111"
$ws.Range("E47").Value = "  -4.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.595"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.393"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05474"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4477"
$ws.Range("E51").Value = "  -1.57%  "
